# Regenerate save_data: update column G ("K", formerly Strike#) with newly
# calculated strike-count values (s_vals) for each row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 2
    7  = 2
    8  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 2
    16 = 2
    17 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 3
    23 = 0
    24 = 2
    25 = 2
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 1
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 3
    36 = 1
    37 = 2
    38 = 1
    40 = 2
    43 = 3
    44 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
